$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.531.04'
$ws.Range("E2").Value = '  -0.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.567.93'
$ws.Range("E3").Value = '  +3.54%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '599.18'
$ws.Range("E5").Value = '  +0.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.59'
$ws.Range("E6").Value = '  +2.82%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.568.67'
$ws.Range("E7").Value = '  +3.56%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("E9").Value = '  +0.94%  '

$ws.Range("E10").Value = '  +3.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.15'
$ws.Range("E11").Value = '  -5.17%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.394'
$ws.Range("E12").Value = '  +3.63%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.175.40'
$ws.Range("E13").Value = '  +3.46%  '

$ws.Range("E14").Value = '  +4.32%  '

$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.14'
$ws.Range("E15").Value = '  +2.29%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.567.80'
$ws.Range("E16").Value = '  +2.64%  '

$ws.Range("E17").Value = '  +1.50%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '65.470.62'
$ws.Range("E18").Value = '  -0.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.27'
$ws.Range("E19").Value = '  +3.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.87'
$ws.Range("E20").Value = '  +1.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.32'
$ws.Range("E21").Value = '  +4.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '396.80'
$ws.Range("E22").Value = '  +0.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.572'
$ws.Range("E23").Value = '  +4.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.712.55'
$ws.Range("E24").Value = '  +3.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '74.70'
$ws.Range("E25").Value = '  +1.97%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("E27").Value = '  +10.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.89'
$ws.Range("E28").Value = '  +8.39%  '

$ws.Range("E29").Value = '  -0.35%  '

$ws.Range("E30").Value = '  +1.09%  '

$ws.Range("E31").Value = '  +1.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.588.15'
$ws.Range("E32").Value = '  +3.87%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.94'
$ws.Range("E33").Value = '  +4.89%  '

$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("E35").Value = '  +1.37%  '

$ws.Range("E36").Value = '  +3.79%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.07'
$ws.Range("E37").Value = '  +2.75%  '

$ws.Range("E38").Value = '  +2.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '168.50'
$ws.Range("E39").Value = '  -2.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.02'
$ws.Range("E40").Value = '  +4.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0804'
$ws.Range("E41").Value = '  +3.53%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.840'
$ws.Range("E42").Value = '  +3.00%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.56'
$ws.Range("E43").Value = '  +15.40%  '

$ws.Range("E44").Value = '  -1.17%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.45'
$ws.Range("E46").Value = '  +0.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.70'
$ws.Range("E47").Value = '  +4.77%  '

$ws.Range("E48").Value = '  +7.87%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.446.23'
$ws.Range("E49").Value = '  +10.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.83'
$ws.Range("E50").Value = '  +4.83%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.16'
$ws.Range("E51").Value = '  +2.44%  '
